# This script normalizes the "Recorded By" column (column G) on the
# "Session Analysis Results" sheet: for each data row, the comma-separated
# list of recorders is reordered so that the literal entry "System" (exact
# case) is moved to the front of the list (preserving the relative order of
# the remaining entries). If "System" is not present in the list, the whole
# list is reversed.
#
# NOTE: PowerShell's default string comparison operators (-eq, -ne,
# -contains, -ceq, ...) behave case-insensitively in this runtime, so we use
# the .Equals() instance method (which is case-sensitive / ordinal) to tell
# "System" apart from "system".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "") { continue }

    $rawParts = $val -split ",\s*"
    $parts = @()
    foreach ($p in $rawParts) { $parts += $p.Trim() }

    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    $newParts = @()
    if ($hasSystem) {
        $systemItems = @()
        $rest = @()
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $systemItems += $p } else { $rest += $p }
        }
        $newParts = $systemItems + $rest
    }
    else {
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
